# MALS-1104 resolve issues found during Demo
#
# The authoring bug: several CDOGS-style merge placeholders were written
# without the dot that separates the repeated-item index from the field
# name, e.g. "{d.Client[i]Phone}" instead of the correct
# "{d.Client[i].Phone}". This edit normalizes every such placeholder to the
# dotted form, and also renames {d.Client[i]TankSerialNo} to
# {d.Client[i].TankSerial}.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Client header block -------------------------------------------------
$ws.Range("B5").Value  = "{d.Client[i].IRMA_Num}"
$ws.Range("B7").Value  = "{d.Client[i].LicenceHolderCompany}"

# --- Contact block ---------------------------------------------------------
$ws.Range("B10").Value = "{d.Client[i].LastnameFirstName}"
$ws.Range("D10").Value = "{d.Client[i].Phone}"

$ws.Range("B11").Value = "{d.Client[i].Address}"
$ws.Range("D11").Value = "{d.Client[i].Fax}"

$ws.Range("B12").Value = "{d.Client[i].City} {d.Client[i].Province}"
$ws.Range("D12").Value = "{d.Client[i].Cell}"

$ws.Range("B13").Value = "{d.Client[i].Postcode}"
$ws.Range("D13").Value = "{d.Client[i].Email}"

# --- Issue date --------------------------------------------------------
$ws.Range("B15").Value = "{d.Client[i].IssueDate}"

# --- Site / tank block ---------------------------------------------------
$ws.Range("B17").Value = "{d.Client[i].SiteAddress}"
$ws.Range("B18").Value = "{d.Client[i].SiteCity} {d.Client[i].SiteProvince}"
$ws.Range("D18").Value = "{d.Client[i].TankCompany}"
$ws.Range("E18").Value = "{d.Client[i].TankModel}"
$ws.Range("F18").Value = "{d.Client[i].TankSerial}"
$ws.Range("G18").Value = "{d.Client[i].TankCapacity}"

# --- Last inspection summary ---------------------------------------------
$ws.Range("B21").Value = "{d.Client[i].LastInspectionDate}"
$ws.Range("D21").Value = "{d.Client[i].LastInspector}"

# --- Inspection history table --------------------------------------------
$ws.Range("A25").Value = "{d.Client[i].Insp[i].Date}"
$ws.Range("B25").Value = "{d.Client[i].Insp[i].IH}"
$ws.Range("C25").Value = "{d.Client[i].Insp[i].SCC}"
$ws.Range("D25").Value = "{d.Client[i].Insp[i].IBC}"
$ws.Range("E25").Value = "{d.Client[i].Insp[i].CRY}"
$ws.Range("A26").Value = "{d.Client[i].Insp[i+1].Date}"

# --- Report average row ---------------------------------------------------
$ws.Range("B28").Value = "{d.Client[i].Avg_IH}"
$ws.Range("C28").Value = "{d.Client[i].Avg_SCC}"
$ws.Range("D28").Value = "{d.Client[i].Avg_IBC}"
$ws.Range("E28").Value = "{d.Client[i].Avg_CRY}"

# --- View state: scroll position & selection ------------------------------
$ws.Range("F17").Select()
